$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(16, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2301", 40000, 1200000),
    @(17, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2212", 48000, 1200000),
    @(18, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2211", 48000, 1200000),
    @(19, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2210", 48000, 1200000),
    @(20, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2209", 48000, 1200000),
    @(21, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2208", 48000, 1200000),
    @(22, "1128048569", "KELLY JOHANNA INFANTE GUZMAN", "2207", 40000, 1200000),
    @(23, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2301", 38867, 1166000),
    @(24, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2212", 46640, 1166000),
    @(25, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2211", 46640, 1166000),
    @(26, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2210", 46640, 1166000),
    @(27, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2209", 46640, 1166000),
    @(28, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2208", 46640, 1166000),
    @(29, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2207", 46640, 1166000),
    @(30, "33221631", "NEYLA ESTHER GUTIERREZ TOSCANO", "2206", 35758, 1166000),
    @(31, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2301", 40333, 1210000),
    @(32, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2212", 48400, 1210000),
    @(33, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2211", 48400, 1210000),
    @(34, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2210", 48400, 1210000),
    @(35, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2209", 48400, 1210000),
    @(36, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2208", 48400, 1210000),
    @(37, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2207", 48400, 1210000),
    @(38, "9297191", "MANUEL DEL CRISTO PAJARO LLERENA", "2206", 37106, 1210000),
    @(39, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2301", 33333, 1000000),
    @(40, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2212", 40000, 1000000),
    @(41, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2211", 40000, 1000000),
    @(42, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2210", 40000, 1000000),
    @(43, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2209", 40000, 1000000),
    @(44, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2208", 40000, 1000000),
    @(45, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2207", 40000, 1000000),
    @(46, "9158605", "ROBINSON ENRIQUE RIBERO HERNANDEZ", "2206", 30666, 1000000),
    @(47, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2301", 40333, 1210000),
    @(48, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2212", 48400, 1210000),
    @(49, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2211", 48400, 1210000),
    @(50, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2210", 48400, 1210000),
    @(51, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2209", 48400, 1210000),
    @(52, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2208", 48400, 1210000),
    @(53, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2207", 48400, 1210000),
    @(54, "1066718197", "JORGE EMILIO PRADO JARAMILLO", "2206", 37106, 1210000),
    @(55, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2301", 33333, 1000000),
    @(56, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2212", 40000, 1000000),
    @(57, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2211", 40000, 1000000),
    @(58, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2210", 40000, 1000000),
    @(59, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2209", 40000, 1000000),
    @(60, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2208", 40000, 1000000),
    @(61, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2207", 40000, 1000000),
    @(62, "1073994971", "VICENTE MANUEL LOZANO FAJARDO", "2206", 30666, 1000000),
    @(63, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2301", 35333, 1060000),
    @(64, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2212", 42400, 1060000),
    @(65, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2211", 42400, 1060000),
    @(66, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2210", 42400, 1060000),
    @(67, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2209", 42400, 1060000),
    @(68, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2208", 42400, 1060000),
    @(69, "1104869717", "NESTOR ADOLFO SANTOS CASTILLO", "2207", 42400, 1060000),
    @(70, "11063562", "SONY RAFAEL RIVERO CORDERO", "2206", 37106, 1160000),
    @(71, "92261028", "MARIO NEL SANTOS CASTILLO", "2301", 33333, 1000000),
    @(72, "92261028", "MARIO NEL SANTOS CASTILLO", "2212", 40000, 1000000),
    @(73, "92261028", "MARIO NEL SANTOS CASTILLO", "2211", 40000, 1000000),
    @(74, "92261028", "MARIO NEL SANTOS CASTILLO", "2210", 40000, 1000000),
    @(75, "92261028", "MARIO NEL SANTOS CASTILLO", "2209", 40000, 1000000),
    @(76, "92261028", "MARIO NEL SANTOS CASTILLO", "2208", 40000, 1000000),
    @(77, "92261028", "MARIO NEL SANTOS CASTILLO", "2207", 40000, 1000000),
    @(78, "92261028", "MARIO NEL SANTOS CASTILLO", "2206", 30666, 1000000),
    @(79, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2301", 34667, 1040000),
    @(80, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2212", 41600, 1040000),
    @(81, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2211", 41600, 1040000),
    @(82, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2210", 41600, 1040000),
    @(83, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2209", 41600, 1040000),
    @(84, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2208", 41600, 1040000),
    @(85, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2207", 41600, 1040000),
    @(86, "1104469717", "NESTOR ADOLFO SANTOS CASTILLO", "2206", 31894, 1040000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[5]   # G: Salario Basico
}
